# Update "想去人数" (interest count) figures across sheets to match the
# latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F17").Value = 96
$ws1.Range("F22").Value = 3997
$ws1.Range("F30").Value = 2903
$ws1.Range("F31").Value = 2080
$ws1.Range("F36").Value = 4014
$ws1.Range("F44").Value = 1541
$ws1.Range("F45").Value = 257
$ws1.Range("F48").Value = 685

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 513

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 153

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 153
$ws4.Range("F18").Value = 96
$ws4.Range("F24").Value = 3997
$ws4.Range("F32").Value = 2903
$ws4.Range("F33").Value = 2080
$ws4.Range("F38").Value = 4014
$ws4.Range("F45").Value = 1541
$ws4.Range("F46").Value = 257
$ws4.Range("F48").Value = 685
